$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Crlf1"
$row2[0,2] = "Cntfr"
$row2[0,3] = "ECs"
$row2[0,4] = 1
$row2[0,5] = 0.3333333333333333
$row2[0,6] = 0.082925
$row2[0,7] = 0.248775
$row2[0,8] = 0.003066867285585202
$row2[0,9] = 0.003066867285585203
$row2[0,10] = 1
$row2[0,11] = 0.3333333333333333
$row2[0,12] = 0.02507166666666667
$row2[0,13] = 0.075215
$row2[0,14] = 0.001520672841173258
$row2[0,15] = 0.001520672841173258
$row2[0,16] = 0.002079067958333334
$row2[0,17] = 0.018711611625
$row2[0,18] = [double]"4.663701788672166E-06"
$row2[0,19] = [double]"4.663701788672168E-06"
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Crlf1"
$row3[0,2] = "Cntfr"
$row3[0,3] = "FAPs"
$row3[0,4] = 1
$row3[0,5] = 0.3333333333333333
$row3[0,6] = 0.082925
$row3[0,7] = 0.248775
$row3[0,8] = 0.003066867285585202
$row3[0,9] = 0.003066867285585203
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 13.42247233333333
$row3[0,13] = 40.267417
$row3[0,14] = 0.8141137727328104
$row3[0,15] = 0.8141137727328106
$row3[0,16] = 1.113058518241667
$row3[0,17] = 10.017526664175
$row3[0,18] = 0.002496778896338603
$row3[0,19] = 0.002496778896338603
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "Crlf1"
$row4[0,2] = "Cntfr"
$row4[0,3] = "sCs"
$row4[0,4] = 1
$row4[0,5] = 0.3333333333333333
$row4[0,6] = 0.082925
$row4[0,7] = 0.248775
$row4[0,8] = 0.003066867285585202
$row4[0,9] = 0.003066867285585203
$row4[0,10] = 3
$row4[0,11] = 1
$row4[0,12] = 3.039675333333333
$row4[0,13] = 9.119026
$row4[0,14] = 0.1843655544260162
$row4[0,15] = 0.1843655544260162
$row4[0,16] = 0.2520650770166666
$row4[0,17] = 2.26858569315
$row4[0,18] = 0.0005654246874579271
$row4[0,19] = 0.0005654246874579273
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = "FAPs"
$row5[0,1] = "Crlf1"
$row5[0,2] = "Cntfr"
$row5[0,3] = "ECs"
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 17.662076
$row5[0,7] = 52.986228
$row5[0,8] = 0.653207634367435
$row5[0,9] = 0.653207634367435
$row5[0,10] = 1
$row5[0,11] = 0.3333333333333333
$row5[0,12] = 0.02507166666666667
$row5[0,13] = 0.075215
$row5[0,14] = 0.001520672841173258
$row5[0,15] = 0.001520672841173258
$row5[0,16] = 0.4428176821133333
$row5[0,17] = 3.98535913902
$row5[0,18] = 0.00099331510922959
$row5[0,19] = 0.00099331510922959
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = "FAPs"
$row6[0,1] = "Crlf1"
$row6[0,2] = "Cntfr"
$row6[0,3] = "FAPs"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 17.662076
$row6[0,7] = 52.986228
$row6[0,8] = 0.653207634367435
$row6[0,9] = 0.653207634367435
$row6[0,10] = 3
$row6[0,11] = 1
$row6[0,12] = 13.42247233333333
$row6[0,13] = 40.267417
$row6[0,14] = 0.8141137727328104
$row6[0,15] = 0.8141137727328106
$row6[0,16] = 237.0687264592306
$row6[0,17] = 2133.618538133076
$row6[0,18] = 0.5317853315927467
$row6[0,19] = 0.5317853315927469
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = "FAPs"
$row7[0,1] = "Crlf1"
$row7[0,2] = "Cntfr"
$row7[0,3] = "sCs"
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 17.662076
$row7[0,7] = 52.986228
$row7[0,8] = 0.653207634367435
$row7[0,9] = 0.653207634367435
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 3.039675333333333
$row7[0,13] = 9.119026
$row7[0,14] = 0.1843655544260162
$row7[0,15] = 0.1843655544260162
$row7[0,16] = 53.68697675265867
$row7[0,17] = 483.182790773928
$row7[0,18] = 0.1204289876654586
$row7[0,19] = 0.1204289876654586
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = "sCs"
$row8[0,1] = "Crlf1"
$row8[0,2] = "Cntfr"
$row8[0,3] = "ECs"
$row8[0,4] = 3
$row8[0,5] = 1
$row8[0,6] = 9.293991
$row8[0,7] = 27.881973
$row8[0,8] = 0.3437254983469798
$row8[0,9] = 0.3437254983469798
$row8[0,10] = 1
$row8[0,11] = 0.3333333333333333
$row8[0,12] = 0.02507166666666667
$row8[0,13] = 0.075215
$row8[0,14] = 0.001520672841173258
$row8[0,15] = 0.001520672841173258
$row8[0,16] = 0.233015844355
$row8[0,17] = 2.097142599195
$row8[0,18] = 0.0005226940301549957
$row8[0,19] = 0.0005226940301549958
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = "sCs"
$row9[0,1] = "Crlf1"
$row9[0,2] = "Cntfr"
$row9[0,3] = "FAPs"
$row9[0,4] = 3
$row9[0,5] = 1
$row9[0,6] = 9.293991
$row9[0,7] = 27.881973
$row9[0,8] = 0.3437254983469798
$row9[0,9] = 0.3437254983469798
$row9[0,10] = 3
$row9[0,11] = 1
$row9[0,12] = 13.42247233333333
$row9[0,13] = 40.267417
$row9[0,14] = 0.8141137727328104
$row9[0,15] = 0.8141137727328106
$row9[0,16] = 124.748337063749
$row9[0,17] = 1122.735033573741
$row9[0,18] = 0.2798316622437251
$row9[0,19] = 0.2798316622437252
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,20
$row10[0,0] = "sCs"
$row10[0,1] = "Crlf1"
$row10[0,2] = "Cntfr"
$row10[0,3] = "sCs"
$row10[0,4] = 3
$row10[0,5] = 1
$row10[0,6] = 9.293991
$row10[0,7] = 27.881973
$row10[0,8] = 0.3437254983469798
$row10[0,9] = 0.3437254983469798
$row10[0,10] = 3
$row10[0,11] = 1
$row10[0,12] = 3.039675333333333
$row10[0,13] = 9.119026
$row10[0,14] = 0.1843655544260162
$row10[0,15] = 0.1843655544260162
$row10[0,16] = 28.250715190922
$row10[0,17] = 254.256436718298
$row10[0,18] = 0.06337114207309964
$row10[0,19] = 0.06337114207309964
$ws.Range("A10:T10").Value = $row10
